$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update footer timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 01:22"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 585469
$ws.Range("C4").Value = 25169
$ws.Range("E4").Value = 525672
$ws.Range("G4").Value = 1487
$ws.Range("H4").Value = 23592

# --- Row 8: Alemania ---
$ws.Range("B8").Value = 130072
$ws.Range("C8").Value = 2218
$ws.Range("E8").Value = 62578
$ws.Range("G8").Value = 172
$ws.Range("H8").Value = 3194

# --- Row 16: Canada ---
$ws.Range("B16").Value = 25680
$ws.Range("C16").Value = 1297
$ws.Range("E16").Value = 17144
$ws.Range("G16").Value = 63
$ws.Range("H16").Value = 780

# --- Row 26: Peru (only Muertes hoy changes) ---
$ws.Range("F26").Value = 143

# --- Japon moves up, inserted right after Peru; Ecuador & Chile shift down one row ---
# Row 27 becomes Japon with fresh data
$ws.Range("A27").Value = "Japon"
$ws.Range("B27").Value = 7618
$ws.Range("C27").Value = 248
$ws.Range("D27").Value = 799
$ws.Range("E27").Value = 6676
$ws.Range("F27").Value = 135
$ws.Range("G27").Value = 20
$ws.Range("H27").Value = 143

# Row 28 becomes Ecuador (old Ecuador data)
$ws.Range("A28").Value = "Ecuador"
$ws.Range("B28").Value = 7529
$ws.Range("C28").Value = 63
$ws.Range("D28").Value = 597
$ws.Range("E28").Value = 6577
$ws.Range("F28").Value = 121
$ws.Range("G28").Value = 22
$ws.Range("H28").Value = 355

# Row 29 becomes Chile (old Chile data)
$ws.Range("A29").Value = "Chile"
$ws.Range("B29").Value = 7525
$ws.Range("C29").Value = 312
$ws.Range("D29").Value = 2367
$ws.Range("E29").Value = 5076
$ws.Range("F29").Value = 387
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 82

# --- Row 35: Chequia ---
$ws.Range("B35").Value = 6059
$ws.Range("C35").Value = 68
$ws.Range("E35").Value = 5397
$ws.Range("F35").Value = 100

# --- Row 78 ---
$ws.Range("B78").Value = 848
$ws.Range("C78").Value = 28
$ws.Range("E78").Value = 738

# --- Row 94 ---
$ws.Range("B94").Value = 548
$ws.Range("C94").Value = 19
$ws.Range("D94").Value = 86
$ws.Range("E94").Value = 449
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 13

# --- Row 96 ---
$ws.Range("B96").Value = 483
$ws.Range("D96").Value = 248
$ws.Range("E96").Value = 227

# --- Libia moves up, inserted right after Sudan; Siria shifts down one row ---
# Row 165 becomes Libia with fresh data
$ws.Range("A165").Value = "Libia"
$ws.Range("B165").Value = 26
$ws.Range("C165").Value = 1
$ws.Range("D165").Value = 9
$ws.Range("E165").Value = 16
$ws.Range("H165").Value = 1

# Row 166 becomes Siria (old Siria data)
$ws.Range("A166").Value = "Siria"
$ws.Range("B166").Value = 25
$ws.Range("C166").Value = 0
$ws.Range("D166").Value = 5
$ws.Range("E166").Value = 18
$ws.Range("H166").Value = 2

# --- Islas Turcas y Caicos moves up, inserted right after Sierra Leona; Cabo Verde & Surinam shift down one row ---
# Row 196 becomes Islas Turcas y Caicos with fresh data
$ws.Range("A196").Value = "Islas Turcas y Caicos"
$ws.Range("B196").Value = 10
$ws.Range("C196").Value = 1
$ws.Range("D196").Value = 0
$ws.Range("E196").Value = 9
$ws.Range("H196").Value = 1

# Row 197 becomes Cabo Verde (old Cabo Verde data)
$ws.Range("A197").Value = "Cabo Verde"
$ws.Range("B197").Value = 10
$ws.Range("C197").Value = 2
$ws.Range("D197").Value = 1
$ws.Range("E197").Value = 8
$ws.Range("H197").Value = 1

# Row 198 becomes Surinam (old Surinam data)
$ws.Range("A198").Value = "Surinam"
$ws.Range("B198").Value = 10
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 6
$ws.Range("E198").Value = 3
$ws.Range("H198").Value = 1
